$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-12 Friday" "2024-07-13 Saturday"

Replace-Text "559÷2=279, 1" "821÷9=91, 2"
Replace-Text "564÷7=80, 4" "416÷8=52, 0"
Replace-Text "726÷6=121, 0" "760÷3=253, 1"
Replace-Text "407÷9=45, 2" "356÷8=44, 4"
Replace-Text "896÷8=112, 0" "665÷9=73, 8"

Replace-Text "811÷9=90, 1" "553÷9=61, 4"
Replace-Text "572÷7=81, 5" "901÷5=180, 1"
Replace-Text "226÷3=75, 1" "598÷7=85, 3"
Replace-Text "635÷6=105, 5" "976÷9=108, 4"
Replace-Text "768÷6=128, 0" "909÷6=151, 3"

Replace-Text "297÷2=148, 1" "686÷8=85, 6"
Replace-Text "460÷8=57, 4" "223÷7=31, 6"
Replace-Text "740÷5=148, 0" "637÷6=106, 1"
Replace-Text "595÷9=66, 1" "811÷3=270, 1"
Replace-Text "561÷7=80, 1" "172÷8=21, 4"

Replace-Text "504÷6=84, 0" "407÷3=135, 2"
Replace-Text "189÷9=21, 0" "412÷4=103, 0"
Replace-Text "875÷2=437, 1" "752÷6=125, 2"
Replace-Text "604÷8=75, 4" "349÷9=38, 7"
Replace-Text "532÷4=133, 0" "585÷7=83, 4"

Replace-Text "277÷9=30, 7" "218÷6=36, 2"
Replace-Text "294÷9=32, 6" "706÷9=78, 4"
Replace-Text "728÷6=121, 2" "542÷9=60, 2"
Replace-Text "427÷3=142, 1" "438÷8=54, 6"
Replace-Text "219÷4=54, 3" "850÷5=170, 0"
